$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

Set-TextValue 'D2' '21.745.18'
Set-TextValue 'E2' '  -1.19%  '

Set-TextValue 'D3' '1.540.41'
Set-TextValue 'E3' '  -0.68%  '

Set-TextValue 'D4' '0.9974'
Set-TextValue 'E4' '  -0.39%  '

Set-TextValue 'D5' '0.9985'
Set-TextValue 'E5' '  -0.23%  '

Set-TextValue 'D6' '290.39'
Set-TextValue 'E6' '  +0.96%  '

Set-TextValue 'D7' '0.3958'
Set-TextValue 'E7' '  +0.67%  '

Set-TextValue 'D8' '0.3202'
Set-TextValue 'E8' '  +0.18%  '

Set-TextValue 'D9' '42.68'
Set-TextValue 'E9' '  +1.27%  '

Set-TextValue 'D10' '0.07193'
Set-TextValue 'E10' '  -1.01%  '

Set-TextValue 'D11' '1.088'
Set-TextValue 'E11' '  -0.17%  '

Set-TextValue 'D12' '0.9999'
Set-TextValue 'E12' '  -0.15%  '

Set-TextValue 'D13' '5.755'
Set-TextValue 'E13' '  +2.72%  '

Set-TextValue 'D14' '18.45'
Set-TextValue 'E14' '  -2.20%  '

Set-TextValue 'D15' '6.656'
Set-TextValue 'E15' '  +0.13%  '

Set-TextValue 'D16' '1.540.16'
Set-TextValue 'E16' '  -0.68%  '

Set-TextValue 'D17' '0.00001101'
Set-TextValue 'E17' '  -1.34%  '

Set-TextValue 'D18' '0.06604'
Set-TextValue 'E18' '  +0.40%  '

Set-TextValue 'D19' '84.55'
Set-TextValue 'E19' '  +0.95%  '

Set-TextValue 'D20' '0.9977'
Set-TextValue 'E20' '  -0.34%  '

Set-TextValue 'D21' '6.165'
Set-TextValue 'E21' '  -1.91%  '

Set-TextValue 'D22' '15.62'
Set-TextValue 'E22' '  -0.43%  '

Set-TextValue 'E23' '  -3.20%  '

Set-TextValue 'D24' '2.378'
Set-TextValue 'E24' '  +1.34%  '

Set-TextValue 'D25' '21.694.91'
Set-TextValue 'E25' '  -1.43%  '

Set-TextValue 'D26' '2.408'
Set-TextValue 'E26' '  -0.62%  '

Set-TextValue 'D27' '151.34'
Set-TextValue 'E27' '  +2.81%  '

Set-TextValue 'D28' '18.50'
Set-TextValue 'E28' '  -0.41%  '

Set-TextValue 'D29' '4.874'
Set-TextValue 'E29' '  +0.83%  '

Set-TextValue 'D30' '1.715.09'
Set-TextValue 'E30' '  -0.54%  '

Set-TextValue 'D31' '117.89'
Set-TextValue 'E31' '  -0.95%  '

Set-TextValue 'D32' '6.141'
Set-TextValue 'E32' '  +8.48%  '

Set-TextValue 'D33' '0.9858'
Set-TextValue 'E33' '  -6.93%  '

Set-TextValue 'D34' '0.08150'
Set-TextValue 'E34' '  -1.96%  '

Set-TextValue 'D35' '8.617'
Set-TextValue 'E35' '  -5.86%  '

Set-TextValue 'D36' '5.217'
Set-TextValue 'E36' '  +2.57%  '

Set-TextValue 'D37' '0.02250'
Set-TextValue 'E37' '  -0.30%  '

Set-TextValue 'E38' '  -6.34%  '

Set-TextValue 'D39' '0.06010'
Set-TextValue 'E39' '  -2.20%  '

Set-TextValue 'D40' '11.35'
Set-TextValue 'E40' '  +7.75%  '

Set-TextValue 'D41' '0.2056'
Set-TextValue 'E41' '  -0.19%  '

Set-TextValue 'D42' '1.192'
Set-TextValue 'E42' '  -1.64%  '

Set-TextValue 'D43' '0.9973'
Set-TextValue 'E43' '  -0.34%  '

Set-TextValue 'D44' '0.5854'
Set-TextValue 'E44' '  +1.19%  '

Set-TextValue 'D45' '13.19'
Set-TextValue 'E45' '  +0.25%  '

Set-TextValue 'D46' '3.729'
Set-TextValue 'E46' '  +0.53%  '

Set-TextValue 'D47' '0.5615'
Set-TextValue 'E47' '  +1.43%  '

Set-TextValue 'B48' 'EOS'
Set-TextValue 'C48' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue 'D48' '1.169'
Set-TextValue 'E48' '  +2.93%  '

Set-TextValue 'B49' 'NEARProtocol'
Set-TextValue 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '1.900'
Set-TextValue 'E49' '  +0.38%  '

Set-TextValue 'D50' '117.20'
Set-TextValue 'E50' '  -0.17%  '

Set-TextValue 'E51' '  -1.14%  '
